$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NumCell($ws, $ref, $val) {
    $ws.Range($ref).Value = $val
}

function Set-TextCell($ws, $ref, $text, $styleSrc) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $text
    $ws.Range($styleSrc).Copy()
    $ws.Range($ref).PasteSpecial(-4122)
}

function Set-NumCellWithStyle($ws, $ref, $val, $styleSrc) {
    $ws.Range($styleSrc).Copy()
    $ws.Range($ref).PasteSpecial(-4122)
    $ws.Range($ref).Value = $val
}

# ---- Header text updates (Volume/Number, Date range) ----
$ws.Range("A8").Characters(21, 1).Text = "7"
$c9 = $ws.Range("C9")
$c9.Characters(27, 8).Text = "2/13/2023"
$c9.Characters(47, 9).Text = "2/19/2023"

# ---- Row 16-30 numeric/text data refresh ----
Set-NumCell $ws "C16" 4
Set-NumCell $ws "D16" 3
Set-NumCell $ws "E16" 33.333333333333
Set-NumCell $ws "F16" 8
Set-NumCell $ws "G16" 9
Set-NumCell $ws "H16" -11.111111111111
Set-NumCell $ws "I16" 15
Set-NumCell $ws "J16" 18
Set-NumCell $ws "K16" -16.666666666666
Set-NumCell $ws "L16" 200
Set-NumCell $ws "M16" -16.666666666666
Set-NumCell $ws "N16" -86.111111111111
Set-NumCell $ws "C17" 2
Set-NumCell $ws "E17" -50
Set-NumCell $ws "F17" 11
Set-NumCell $ws "G17" 12
Set-NumCell $ws "H17" -8.333333333333
Set-NumCell $ws "I17" 23
Set-NumCell $ws "J17" 20
Set-NumCell $ws "K17" 15
Set-NumCell $ws "L17" 0
Set-NumCell $ws "M17" 76.923076923076
Set-NumCell $ws "N17" 53.333333333333
Set-NumCell $ws "C18" 3
Set-NumCell $ws "D18" 5
Set-NumCell $ws "E18" -40
Set-NumCell $ws "F18" 8
Set-NumCell $ws "G18" 13
Set-NumCell $ws "H18" -38.461538461538
Set-NumCell $ws "I18" 12
Set-NumCell $ws "J18" 22
Set-NumCell $ws "K18" -45.454545454545
Set-NumCell $ws "L18" 0
Set-NumCell $ws "M18" 71.428571428571
Set-NumCell $ws "N18" -75
Set-NumCell $ws "C19" 7
Set-NumCell $ws "D19" 14
Set-NumCell $ws "E19" -50
Set-NumCell $ws "F19" 29
Set-NumCell $ws "G19" 61
Set-NumCell $ws "H19" -52.459016393442
Set-NumCell $ws "I19" 67
Set-NumCell $ws "J19" 99
Set-NumCell $ws "K19" -32.323232323232
Set-NumCell $ws "L19" 52.272727272727
Set-NumCell $ws "M19" 168
Set-NumCell $ws "N19" 31.372549019607
Set-NumCellWithStyle $ws "D20" 1 "F15"
Set-NumCellWithStyle $ws "E20" 0 "H15"
Set-NumCell $ws "I20" 7
Set-NumCell $ws "J20" 6
Set-NumCell $ws "K20" 16.666666666666
Set-NumCell $ws "L20" -12.5
Set-NumCell $ws "M20" -12.5
Set-NumCell $ws "N20" -86.538461538461
Set-NumCell $ws "C21" 17
Set-NumCell $ws "D21" 27
Set-NumCell $ws "E21" -37.037037037037
Set-NumCell $ws "F21" 61
Set-NumCell $ws "G21" 97
Set-NumCell $ws "H21" -37.113402061855
Set-NumCell $ws "I21" 126
Set-NumCell $ws "J21" 168
Set-NumCell $ws "K21" -25
Set-NumCell $ws "L21" 35.483870967741
Set-NumCell $ws "M21" 72.602739726027
Set-NumCell $ws "N21" -54.512635379061
Set-TextCell $ws "D22" "0" "C15"
Set-TextCell $ws "E22" "***.*" "E15"
Set-NumCell $ws "G22" 2
Set-NumCellWithStyle $ws "L22" -50 "H15"
Set-NumCell $ws "C23" 6
Set-NumCell $ws "E23" 200
Set-NumCell $ws "F23" 13
Set-NumCell $ws "G23" 8
Set-NumCell $ws "H23" 62.5
Set-NumCell $ws "I23" 20
Set-NumCell $ws "J23" 15
Set-NumCell $ws "K23" 33.333333333333
Set-NumCell $ws "L23" -4.761904761904
Set-NumCell $ws "M23" 11.111111111111
Set-NumCell $ws "C24" 21
Set-NumCell $ws "D24" 50
Set-NumCell $ws "E24" -58
Set-NumCell $ws "F24" 77
Set-NumCell $ws "G24" 206
Set-NumCell $ws "H24" -62.621359223301
Set-NumCell $ws "I24" 148
Set-NumCell $ws "J24" 333
Set-NumCell $ws "K24" -55.555555555555
Set-NumCell $ws "L24" 27.586206896551
Set-NumCell $ws "M24" 78.313253012048
Set-NumCell $ws "F25" 35
Set-NumCell $ws "G25" 27
Set-NumCell $ws "H25" 29.629629629629
Set-NumCell $ws "I25" 52
Set-NumCell $ws "J25" 43
Set-NumCell $ws "K25" 20.930232558139
Set-NumCell $ws "L25" 44.444444444444
Set-TextCell $ws "D27" "0" "C15"
Set-TextCell $ws "E27" "***.*" "E15"
Set-NumCell $ws "L27" -50
Set-TextCell $ws "C28" "0" "C15"
Set-NumCellWithStyle $ws "L28" 100 "H15"
Set-TextCell $ws "C29" "0" "C15"
Set-NumCellWithStyle $ws "L29" 100 "H15"
Set-TextCell $ws "G30" "0" "C15"
Set-TextCell $ws "H30" "***.*" "E15"

Write-Output "done"